$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for all data rows (2..305) to 45180
$ws.Range("C2:C305").Value = 45180

# 2. Insert a new row at position 7. This shifts the existing row 7
#    ("A 64630-2019") down to row 8, and old row 8 ("A 39885-2023") down to row 9.
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).RowHeight = 15

# 3. Populate the newly inserted row 7 with the (updated) "A 39885-2023" record.
$ws.Cells.Item(7,1).Value = "A 39885-2023"
$ws.Cells.Item(7,2).Value = 45168
$ws.Cells.Item(7,3).Value = 45180
$ws.Cells.Item(7,4).Value = "GÄVLEBORGS LÄN"
$ws.Cells.Item(7,5).Value = "SÖDERHAMN"
$ws.Cells.Item(7,7).Value = 2.9
$ws.Cells.Item(7,8).Value = 0
$ws.Cells.Item(7,9).Value = 2
$ws.Cells.Item(7,10).Value = 1
$ws.Cells.Item(7,11).Value = 0
$ws.Cells.Item(7,12).Value = 0
$ws.Cells.Item(7,13).Value = 0
$ws.Cells.Item(7,14).Value = 0
$ws.Cells.Item(7,15).Value = 1
$ws.Cells.Item(7,16).Value = 0
$ws.Cells.Item(7,17).Value = 3
$ws.Cells.Item(7,18).Value = "Hapalopilus aurantiacus`r`nBlodticka`r`nSotriska"
$ws.Cells.Item(7,19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SODERHAMN/artfynd/A 39885-2023.xlsx")'
$ws.Cells.Item(7,20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SODERHAMN/kartor/A 39885-2023.png")'
$ws.Cells.Item(7,22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SODERHAMN/klagomål/A 39885-2023.docx")'
$ws.Cells.Item(7,23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SODERHAMN/klagomålsmail/A 39885-2023.docx")'
$ws.Cells.Item(7,24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SODERHAMN/tillsyn/A 39885-2023.docx")'
$ws.Cells.Item(7,25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_SODERHAMN/tillsynsmail/A 39885-2023.docx")'

# Row 7 should use the date number format for B/C, same as other rows,
# and wrap-text for the species column R - matching the formatting used
# throughout the rest of the sheet's data rows.
$ws.Cells.Item(7,2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(7,3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(7,18).WrapText = $true

# 4. The old "A 39885-2023" row has now been duplicated at row 9 (it shifted
#    down twice: once because of the insert). Remove that duplicate so the
#    original "A 63515-2018" record becomes row 9 again, matching the target.
$ws.Rows.Item(9).Delete()
